$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised values for existing rows 419-421
$ws.Range("C419:F419").Value = 5245350000000
$ws.Range("C420:F420").Value = 5322265000000
$ws.Range("C421:F421").Value = 5421638000000

# Append new rows 422-424 (same style as the last existing data row)
$newRows = @(
    @(45108.41666666666, "ECONOMICS:BRM2", 5501072000000),
    @(45139.41666666666, "ECONOMICS:BRM2", 5591097000000),
    @(45170.41666666666, "ECONOMICS:BRM2", 5656835000000)
)

$r = 422
foreach ($row in $newRows) {
    $dt = $row[0]
    $sym = $row[1]
    $val = $row[2]

    # Carry the date-column cell formatting (style index) down from the row above
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $dt
    $ws.Cells.Item($r, 2).Value = $sym
    $ws.Range("C" + $r + ":F" + $r).Value = $val
    $ws.Cells.Item($r, 7).Value = 0

    $r = $r + 1
}
